$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "encoded" column for 1 Chronicles / 2 Chronicles switches from a plain 0
# to the "*" placeholder marker (used elsewhere for not-yet-geocoded rows).
$ws.Range("F14").Value = "*"
$ws.Range("F15").Value = "*"

# Song of Solomon ("Cantares") finishes geocoding: encoded flag becomes 1
# and latitude/longitude coordinates are filled in.
$ws.Range("F23").Value = 1
$ws.Range("K23").Value = 32.2361352458334
$ws.Range("L23").Value = 35.5485861263889

# Match the row height used by the rest of the summary rows.
$ws.Rows.Item(68).RowHeight = 13.8

# Move the saved selection to the last totals cell.
$ws.Range("I71").Select() | Out-Null
